$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D5").Value = "2계 선형 미분방정식의 해법"
$ws.Range("E5").Value = "https://angeloyeo.github.io/2021/05/27/second_order_ODE.html"

$ws.Range("D6").Value = "Deep Learning for Time Series Forecasting (kaggle 코드 리뷰)"
$ws.Range("E6").Value = "https://leedakyeong.tistory.com/entry/Deep-Learning-for-Time-Series-Forecasting-kaggle"

$ws.Range("D21").Value = "[kaldi] allocate_egs.py 무한루프 현상"
$ws.Range("E21").Value = "https://ms-review.tistory.com/14"
